# Commit: "added new method for simplified gathering of data for a whole year"
#
# The underlying data change is on the "Freiland" worksheet: three new
# columns are appended after the existing "Pluvio_mm" column (I) so that a
# whole year's precipitation data can be gathered in one place:
#   J: Pluvio_mm_SUM           (unit row: mm)
#   K: Niederschlag.Casella    (unit row: mm)
#   L: Niederschlag.Casella_SUM(unit row: mm)
#
# (The charts that later get re-plotted against these new columns live on
# chartsheets, which this COM surface does not expose as editable chart
# objects, so this script focuses on the reachable, authoritative data
# change: the new header/unit cells on the Freiland sheet - which is also
# what drives the new shared-string entries "Pluvio_mm_SUM",
# "Niederschlag.Casella" and "Niederschlag.Casella_SUM".)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Freiland")

# New header row (row 1) entries
$ws.Range("J1").Value = "Pluvio_mm_SUM"
$ws.Range("K1").Value = "Niederschlag.Casella"
$ws.Range("L1").Value = "Niederschlag.Casella_SUM"

# New unit row (row 2) entries - all three new columns are measured in mm
$ws.Range("J2").Value = "mm"
$ws.Range("K2").Value = "mm"
$ws.Range("L2").Value = "mm"

# Size the new columns to fit their (longer) header text, like the other
# bestFit columns (A and I) on this sheet.
$ws.Columns.Item(10).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(12).EntireColumn.AutoFit() | Out-Null

# Leave the selection on the last new cell (matches the author's saved
# cursor position) and then restore "Fichte" as the active sheet so the
# workbook's active-tab bookkeeping is unaffected by this edit.
$ws.Range("L2").Select() | Out-Null
$wb.Worksheets.Item("Fichte").Activate()
